# Append the new data row (homework row for 2020-03-12) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 97

# timestamp (numeric)
$ws.Cells.Item($row, 1).Value = 1583971200
# date / id are stored as text in this sheet, not as Excel dates/numbers,
# so force text with a leading apostrophe (classic Excel "treat as text" trick)
$ws.Cells.Item($row, 2).Value = "'2020-03-12"
$ws.Cells.Item($row, 3).Value = "'5293"
$ws.Cells.Item($row, 4).Value = "AME"
# OHLC + volume (numeric)
$ws.Cells.Item($row, 5).Value = 1.53
$ws.Cells.Item($row, 6).Value = 1.59
$ws.Cells.Item($row, 7).Value = 1.51
$ws.Cells.Item($row, 8).Value = 1.53
$ws.Cells.Item($row, 9).Value = 1043700
